$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 2666.6667
$ws.Cells.Item(64, 9).Value = 2600
$ws.Cells.Item(64, 11).Value = 2600
$ws.Cells.Item(64, 13).Value = -2352

$ws.Cells.Item(67, 8).Value = 2666.6667
$ws.Cells.Item(67, 9).Value = 2600
$ws.Cells.Item(67, 11).Value = 2600
$ws.Cells.Item(67, 13).Value = -1742

$ws.Cells.Item(132, 8).Value = 10200.038
$ws.Cells.Item(132, 9).Value = 10741.471
$ws.Cells.Item(132, 10).Value = 9177.333
$ws.Cells.Item(132, 11).Value = 32224.413
$ws.Cells.Item(132, 12).Value = 27531.999
$ws.Cells.Item(132, 13).Value = -29694.413
$ws.Cells.Item(132, 14).Value = -32591.999

$ws.Cells.Item(138, 8).Value = 1734.4507
$ws.Cells.Item(138, 9).Value = 732.63635
$ws.Cells.Item(138, 10).Value = 3367.037
$ws.Cells.Item(138, 11).Value = 2197.90905
$ws.Cells.Item(138, 12).Value = 10101.111
$ws.Cells.Item(138, 13).Value = 2942.09095
$ws.Cells.Item(138, 14).Value = -20381.111

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(38, 8).Value = 7713.857
$ws.Cells.Item(38, 9).Value = 7713.857
$ws.Cells.Item(38, 11).Value = 7713.857
$ws.Cells.Item(38, 13).Value = -7246.857

$ws.Cells.Item(74, 8).Value = 1234.7413
$ws.Cells.Item(74, 9).Value = 1271.3658
$ws.Cells.Item(74, 10).Value = 1146.4117
$ws.Cells.Item(74, 11).Value = 1271.3658
$ws.Cells.Item(74, 12).Value = 1146.4117
$ws.Cells.Item(74, 13).Value = -397.3658
$ws.Cells.Item(74, 14).Value = -2894.4117

$ws.Cells.Item(77, 8).Value = 1234.7413
$ws.Cells.Item(77, 9).Value = 1271.3658
$ws.Cells.Item(77, 10).Value = 1146.4117
$ws.Cells.Item(77, 11).Value = 6356.829
$ws.Cells.Item(77, 12).Value = 5732.058500000001
$ws.Cells.Item(77, 13).Value = -1988.829
$ws.Cells.Item(77, 14).Value = -14468.0585

$ws.Cells.Item(132, 8).Value = 1132354.2
$ws.Cells.Item(132, 9).Value = 1031.1538
$ws.Cells.Item(132, 10).Value = 4526323.5
$ws.Cells.Item(132, 11).Value = 3093.4614
$ws.Cells.Item(132, 12).Value = 13578970.5
$ws.Cells.Item(132, 13).Value = -563.4614000000001
$ws.Cells.Item(132, 14).Value = -13584030.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 171.5
$ws.Cells.Item(19, 9).Value = 135
$ws.Cells.Item(19, 10).Value = 500
$ws.Cells.Item(19, 11).Value = 135
$ws.Cells.Item(19, 12).Value = 500
$ws.Cells.Item(19, 13).Value = 35
$ws.Cells.Item(19, 14).Value = -840

$ws.Cells.Item(24, 8).Value = 171.5
$ws.Cells.Item(24, 9).Value = 135
$ws.Cells.Item(24, 10).Value = 500
$ws.Cells.Item(24, 11).Value = 135
$ws.Cells.Item(24, 12).Value = 500
$ws.Cells.Item(24, 13).Value = 35
$ws.Cells.Item(24, 14).Value = -840

$ws.Cells.Item(32, 8).Value = 5400
$ws.Cells.Item(32, 9).Value = 3440
$ws.Cells.Item(32, 11).Value = 3440
$ws.Cells.Item(32, 13).Value = -3124

$ws.Cells.Item(58, 8).Value = 30303748
$ws.Cells.Item(58, 9).Value = 62500710
$ws.Cells.Item(58, 10).Value = 722.2941
$ws.Cells.Item(58, 11).Value = 62500710
$ws.Cells.Item(58, 12).Value = 722.2941
$ws.Cells.Item(58, 13).Value = -62500507
$ws.Cells.Item(58, 14).Value = -1128.2941

$ws.Cells.Item(62, 8).Value = 4010.5
$ws.Cells.Item(62, 9).Value = 2776
$ws.Cells.Item(62, 10).Value = 5245
$ws.Cells.Item(62, 11).Value = 2776
$ws.Cells.Item(62, 12).Value = 5245
$ws.Cells.Item(62, 13).Value = -2152
$ws.Cells.Item(62, 14).Value = -6493

$ws.Cells.Item(65, 8).Value = 4010.5
$ws.Cells.Item(65, 9).Value = 2776
$ws.Cells.Item(65, 10).Value = 5245
$ws.Cells.Item(65, 11).Value = 13880
$ws.Cells.Item(65, 12).Value = 26225
$ws.Cells.Item(65, 13).Value = -10760
$ws.Cells.Item(65, 14).Value = -32465

$ws.Cells.Item(134, 8).Value = 949.7143
$ws.Cells.Item(134, 9).Value = 1008
$ws.Cells.Item(134, 10).Value = 600
$ws.Cells.Item(134, 11).Value = 3024
$ws.Cells.Item(134, 12).Value = 1800
$ws.Cells.Item(134, 13).Value = -489
$ws.Cells.Item(134, 14).Value = -6870

$ws.Cells.Item(136, 8).Value = 30303748
$ws.Cells.Item(136, 9).Value = 62500710
$ws.Cells.Item(136, 10).Value = 722.2941
$ws.Cells.Item(136, 11).Value = 187502130
$ws.Cells.Item(136, 12).Value = 2166.8823
$ws.Cells.Item(136, 13).Value = -187499580
$ws.Cells.Item(136, 14).Value = -7266.882299999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(103, 8).Value = 2987.9375
$ws.Cells.Item(103, 10).Value = 3729.8333
$ws.Cells.Item(103, 12).Value = 11189.4999
$ws.Cells.Item(103, 14).Value = -12947.4999

$ws.Cells.Item(114, 8).Value = 733.7273
$ws.Cells.Item(114, 9).Value = 248.375
$ws.Cells.Item(114, 10).Value = 1011.0714
$ws.Cells.Item(114, 11).Value = 745.125
$ws.Cells.Item(114, 12).Value = 3033.2142
$ws.Cells.Item(114, 13).Value = 2508.875
$ws.Cells.Item(114, 14).Value = -9541.2142

$ws.Cells.Item(130, 8).Value = 2022.2222
$ws.Cells.Item(130, 10).Value = 2357.1428
$ws.Cells.Item(130, 12).Value = 7071.428400000001
$ws.Cells.Item(130, 14).Value = -17111.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 5209.6
$ws.Cells.Item(132, 9).Value = 2249.3333
$ws.Cells.Item(132, 10).Value = 8344
$ws.Cells.Item(132, 11).Value = 6747.999899999999
$ws.Cells.Item(132, 12).Value = 25032
$ws.Cells.Item(132, 13).Value = -4217.999899999999
$ws.Cells.Item(132, 14).Value = -30092

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 2650
$ws.Cells.Item(32, 9).Value = 2500
$ws.Cells.Item(32, 10).Value = 2800
$ws.Cells.Item(32, 11).Value = 2500
$ws.Cells.Item(32, 12).Value = 2800
$ws.Cells.Item(32, 13).Value = -2183
$ws.Cells.Item(32, 14).Value = -3434

$ws.Cells.Item(69, 8).Value = 37265.332
$ws.Cells.Item(69, 10).Value = 37265.332
$ws.Cells.Item(69, 12).Value = 37265.332
$ws.Cells.Item(69, 14).Value = -38887.332

$ws.Cells.Item(72, 8).Value = 37265.332
$ws.Cells.Item(72, 10).Value = 37265.332
$ws.Cells.Item(72, 12).Value = 111795.996
$ws.Cells.Item(72, 14).Value = -119907.996

$ws.Cells.Item(132, 8).Value = 5829.302
$ws.Cells.Item(132, 9).Value = 1424.4286
$ws.Cells.Item(132, 10).Value = 14394.333
$ws.Cells.Item(132, 11).Value = 4273.2858
$ws.Cells.Item(132, 12).Value = 43182.999
$ws.Cells.Item(132, 13).Value = -1743.2858
$ws.Cells.Item(132, 14).Value = -48242.999

$ws.Cells.Item(136, 8).Value = 32654982
$ws.Cells.Item(136, 9).Value = 5104131.5
$ws.Cells.Item(136, 10).Value = 142858380
$ws.Cells.Item(136, 11).Value = 15312394.5
$ws.Cells.Item(136, 12).Value = 428575140
$ws.Cells.Item(136, 13).Value = -15309844.5
$ws.Cells.Item(136, 14).Value = -428580240

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(28, 8).Value = 1017
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 14).Value = ""

$ws.Cells.Item(136, 8).Value = 13890818
$ws.Cells.Item(136, 9).Value = 22728278
$ws.Cells.Item(136, 10).Value = 3379.2856
$ws.Cells.Item(136, 11).Value = 68184834
$ws.Cells.Item(136, 12).Value = 10137.8568
$ws.Cells.Item(136, 13).Value = -68182284
$ws.Cells.Item(136, 14).Value = -15237.8568
